$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.537.74"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.511.93"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("D7").Value = "3.509.85"
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  +3.44%  "
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.129"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.403"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("D13").Value = "4.100.22"
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000194"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("D16").Value = "3.503.20"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "66.386.48"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.585"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").Value = "3.657.23"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.69%  "
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "3.516.95"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.37%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("E38").Value = "  -5.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "172.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0801"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.850"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("E45").Value = "  -7.58%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -9.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("E49").Value = "  -5.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.893"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.52%  "
